$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "en informatique près de Grenoble.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "près de Grenoble.",
    2
)
